$wb = $excel.ActiveWorkbook

# This script applies numeric corrections to the currentAveragePrice(NQ/HQ)
# and derived LevePrice/LeveProfit columns (H, I, J, K, L, M, N) across several
# worksheets, reflecting refreshed market-board pricing data pulled by the
# scheduled runner. Each block targets one worksheet and writes the updated
# cell values directly (source data, not formulas).

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 629.2143
$ws.Range("I15").Value = 629.2143
$ws.Range("K15").Value = 1887.6429
$ws.Range("M15").Value = -1718.6429
$ws.Range("H43").Value = 4317
$ws.Range("J43").Value = 4152.4165
$ws.Range("L43").Value = 4152.4165
$ws.Range("N43").Value = -4290.4165
$ws.Range("H100").Value = 30366.334
$ws.Range("I100").Value = 38792.85
$ws.Range("K100").Value = 38792.85
$ws.Range("M100").Value = -38251.85
$ws.Range("H129").Value = 1347.85
$ws.Range("I129").Value = 998.4
$ws.Range("K129").Value = 2995.2
$ws.Range("M129").Value = 2004.8
$ws.Range("H131").Value = 9225
$ws.Range("I131").Value = 5950
$ws.Range("K131").Value = 17850
$ws.Range("M131").Value = -12810
$ws.Range("H141").Value = 2269.6
$ws.Range("I141").Value = 2246
$ws.Range("K141").Value = 6738
$ws.Range("M141").Value = -1558

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 10660.467
$ws.Range("I45").Value = 11738
$ws.Range("K45").Value = 11738
$ws.Range("M45").Value = -11361
$ws.Range("H61").Value = 8577.259
$ws.Range("J61").Value = 14278.556
$ws.Range("L61").Value = 14278.556
$ws.Range("N61").Value = -14702.556
$ws.Range("H63").Value = 6500.3335
$ws.Range("J63").Value = 8753
$ws.Range("L63").Value = 8753
$ws.Range("N63").Value = -10125
$ws.Range("H66").Value = 6500.3335
$ws.Range("J66").Value = 8753
$ws.Range("L66").Value = 43765
$ws.Range("N66").Value = -50629
$ws.Range("H74").Value = 14439.115
$ws.Range("I74").Value = 16725.2
$ws.Range("J74").Value = 6818.8335
$ws.Range("K74").Value = 16725.2
$ws.Range("L74").Value = 6818.8335
$ws.Range("M74").Value = -15851.2
$ws.Range("N74").Value = -8566.833500000001
$ws.Range("H77").Value = 14439.115
$ws.Range("I77").Value = 16725.2
$ws.Range("J77").Value = 6818.8335
$ws.Range("K77").Value = 83626
$ws.Range("L77").Value = 34094.1675
$ws.Range("M77").Value = -79258
$ws.Range("N77").Value = -42830.1675
$ws.Range("H132").Value = 11935.143
$ws.Range("I132").Value = 7257.8335
$ws.Range("K132").Value = 21773.5005
$ws.Range("M132").Value = -19243.5005
$ws.Range("H136").Value = 8577.259
$ws.Range("J136").Value = 14278.556
$ws.Range("L136").Value = 42835.66800000001
$ws.Range("N136").Value = -47935.66800000001

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 425.8
$ws.Range("I4").Value = 425.8
$ws.Range("K4").Value = 425.8
$ws.Range("M4").Value = -310.8
$ws.Range("H35").Value = 96995
$ws.Range("J35").Value = 96995
$ws.Range("L35").Value = 96995
$ws.Range("N35").Value = -97615
$ws.Range("H81").Value = 22665.666
$ws.Range("J81").Value = 22665.666
$ws.Range("L81").Value = 22665.666
$ws.Range("N81").Value = -24787.666
$ws.Range("H82").Value = 60509.688
$ws.Range("J82").Value = 78400.086
$ws.Range("L82").Value = 78400.086
$ws.Range("N82").Value = -79166.086
$ws.Range("H84").Value = 22665.666
$ws.Range("J84").Value = 22665.666
$ws.Range("L84").Value = 67996.99800000001
$ws.Range("N84").Value = -78604.99800000001
$ws.Range("H85").Value = 60509.688
$ws.Range("J85").Value = 78400.086
$ws.Range("L85").Value = 78400.086
$ws.Range("N85").Value = -81052.086
$ws.Range("H134").Value = 11636.448
$ws.Range("I134").Value = 5350.2354
$ws.Range("K134").Value = 16050.7062
$ws.Range("M134").Value = -13515.7062

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2843.1228
$ws.Range("I58").Value = 1799.5667
$ws.Range("J58").Value = 4002.6296
$ws.Range("K58").Value = 1799.5667
$ws.Range("L58").Value = 4002.6296
$ws.Range("M58").Value = -1596.5667
$ws.Range("N58").Value = -4408.6296
$ws.Range("H94").Value = 8245.182000000001
$ws.Range("I94").Value = 7659.8
$ws.Range("J94").Value = 8733
$ws.Range("K94").Value = 7659.8
$ws.Range("L94").Value = 8733
$ws.Range("M94").Value = -7208.8
$ws.Range("N94").Value = -9635
$ws.Range("H134").Value = 4174.467
$ws.Range("I134").Value = 3069.111
$ws.Range("K134").Value = 9207.332999999999
$ws.Range("M134").Value = -6672.332999999999
$ws.Range("H136").Value = 2843.1228
$ws.Range("I136").Value = 1799.5667
$ws.Range("J136").Value = 4002.6296
$ws.Range("K136").Value = 5398.7001
$ws.Range("L136").Value = 12007.8888
$ws.Range("M136").Value = -2848.7001
$ws.Range("N136").Value = -17107.8888

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H102").Value = 7225.5
$ws.Range("I102").Value = 7225.5
$ws.Range("K102").Value = 21676.5
$ws.Range("M102").Value = -19242.5
$ws.Range("H103").Value = 419.5
$ws.Range("I103").Value = 419
$ws.Range("K103").Value = 1257
$ws.Range("M103").Value = -378
$ws.Range("H104").Value = 5619.5
$ws.Range("I104").Value = 5299.5557
$ws.Range("K104").Value = 15898.6671
$ws.Range("M104").Value = -13277.6671
$ws.Range("H108").Value = 2819.8
$ws.Range("I108").Value = 2524.75
$ws.Range("K108").Value = 7574.25
$ws.Range("M108").Value = -4694.25
$ws.Range("H109").Value = 4157
$ws.Range("I109").Value = 4157
$ws.Range("K109").Value = 12471
$ws.Range("M109").Value = -11431
$ws.Range("H110").Value = 17539.6
$ws.Range("I110").Value = 17539.6
$ws.Range("K110").Value = 52618.8
$ws.Range("M110").Value = -48528.8
$ws.Range("H111").Value = 999.5
$ws.Range("I111").Value = 999.5
$ws.Range("K111").Value = 2998.5
$ws.Range("M111").Value = 68.5
$ws.Range("H112").Value = 847.5
$ws.Range("J112").Value = 1495
$ws.Range("L112").Value = 4485
$ws.Range("N112").Value = -6701
$ws.Range("H114").Value = 2776.3333
$ws.Range("I114").Value = 800
$ws.Range("J114").Value = 3341
$ws.Range("K114").Value = 2400
$ws.Range("L114").Value = 10023
$ws.Range("M114").Value = 854
$ws.Range("N114").Value = -16531
$ws.Range("H115").Value = 1677.3636
$ws.Range("I115").Value = 1411.3334
$ws.Range("J115").Value = 1996.6
$ws.Range("K115").Value = 4234.0002
$ws.Range("L115").Value = 5989.799999999999
$ws.Range("M115").Value = -3059.0002
$ws.Range("N115").Value = -8339.799999999999

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 17133.818
$ws.Range("I132").Value = 8947.200000000001
$ws.Range("K132").Value = 26841.6
$ws.Range("M132").Value = -24311.6

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 5534.25
$ws.Range("J122").Value = 3568.5
$ws.Range("L122").Value = 10705.5
$ws.Range("N122").Value = -15605.5
